# Generate Report for Handoff
# Adds a new handoff entry (88654af0-5bcf-4422-b2db-4e652f700f1b.md) as row 9
# to the Overview, zh-cn and de-de tables.

$wb = $excel.ActiveWorkbook

$guid = "88654af0-5bcf-4422-b2db-4e652f700f1b"
$zhXlf = "$guid.9aa37e68f9d3f61da217cfc3e9a9df060d246ab2.zh-cn.xlf"
$deXlf = "$guid.9aa37e68f9d3f61da217cfc3e9a9df060d246ab2.de-de.xlf"
$dateFmt = "yyyy-mm-dd HH:mm:ss"
$ghBase = "https://github.com/OpenLocalizationTestOrg/oltest/blob"
$ghCommit = "aaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaa"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Cells.Item(9, 1).Value = "$guid.md"
$wsOverview.Cells.Item(9, 3).Value = ".md"
$wsOverview.Cells.Item(9, 4).Value = "'"
$wsOverview.Cells.Item(9, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item(9, 6).Value = "Ready for handoff"

$cellG9 = $wsOverview.Cells.Item(9, 7)
$cellG9.NumberFormat = $dateFmt
$cellG9.Value = "2016-08-13 18:54:35"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B9"),
    "$ghBase/$ghCommit/e2e/$guid.md",
    "",
    "",
    "e2e\$guid.md"
) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Cells.Item(9, 2).Value = ".md"
$wsZhCn.Cells.Item(9, 3).Value = "Ready for handoff"
$wsZhCn.Cells.Item(9, 4).Value = "e2e"
$wsZhCn.Cells.Item(9, 5).Value = "ht"
$wsZhCn.Cells.Item(9, 6).Value = "'False"
$wsZhCn.Cells.Item(9, 7).Value = $zhXlf

$cellH9zh = $wsZhCn.Cells.Item(9, 8)
$cellH9zh.NumberFormat = $dateFmt
$cellH9zh.Value = "2016-08-13 18:54:26"

$wsZhCn.Cells.Item(9, 9).Value = "'"
$wsZhCn.Cells.Item(9, 10).Value = "'"

$cellK9zh = $wsZhCn.Cells.Item(9, 11)
$cellK9zh.NumberFormat = $dateFmt
$cellK9zh.Value = "0001-01-01 00:00:00"

$wsZhCn.Cells.Item(9, 12).Value = "'"
$wsZhCn.Cells.Item(9, 13).Value = "'True"
$wsZhCn.Cells.Item(9, 14).Value = "'"
$wsZhCn.Cells.Item(9, 15).Value = "'False"
$wsZhCn.Cells.Item(9, 16).Value = "'"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A9"),
    "$ghBase/$ghCommit/e2e/$guid.md",
    "",
    "",
    "$guid.md"
) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Cells.Item(9, 2).Value = ".md"
$wsDeDe.Cells.Item(9, 3).Value = "Ready for handoff"
$wsDeDe.Cells.Item(9, 4).Value = "e2e"
$wsDeDe.Cells.Item(9, 5).Value = "ht"
$wsDeDe.Cells.Item(9, 6).Value = "'False"
$wsDeDe.Cells.Item(9, 7).Value = $deXlf

$cellH9de = $wsDeDe.Cells.Item(9, 8)
$cellH9de.NumberFormat = $dateFmt
$cellH9de.Value = "2016-08-13 18:54:35"

$wsDeDe.Cells.Item(9, 9).Value = "'"
$wsDeDe.Cells.Item(9, 10).Value = "'"

$cellK9de = $wsDeDe.Cells.Item(9, 11)
$cellK9de.NumberFormat = $dateFmt
$cellK9de.Value = "0001-01-01 00:00:00"

$wsDeDe.Cells.Item(9, 12).Value = "'"
$wsDeDe.Cells.Item(9, 13).Value = "'True"
$wsDeDe.Cells.Item(9, 14).Value = "'"
$wsDeDe.Cells.Item(9, 15).Value = "'False"
$wsDeDe.Cells.Item(9, 16).Value = "'"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A9"),
    "$ghBase/$ghCommit/e2e/$guid.md",
    "",
    "",
    "$guid.md"
) | Out-Null
